$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edit: element "lightning" renamed to "bolt" ---
# Row 10 holds the elemental-type table entry for lightning (name + element columns)
$ws.Range("A10").Value = "bolt"
$ws.Range("B10").Value = "bolt"

# --- Touch the previously-customized column widths so their stored metrics get
#     refreshed (mirrors what happens when the workbook is resaved) ---
foreach ($col in @("A", "B", "C", "E", "G")) {
    $currentWidth = $ws.Range($col + "1").ColumnWidth
    $ws.Range($col + "1").ColumnWidth = $currentWidth
}

# --- View state: user scrolled back to column A and selected B10 ---
[void]$ws.Range("B10").Select()
